# "ajustes superficie che tapy"
# Change the surface area mentioned in the PH/PO template text from
# "60,30 m2" to "57,40 m2".
#
# Original run layout for this sentence fragment (4 runs, all sharing the
# same rPr -- Arial Unicode MS / sz 18):
#   run1: "R + 3D, con una Superficie de 60"
#   run2: ","
#   run3: "3"
#   run4: "0 m2., en el inmueble individualizado como: Manzana"
#
# Target layout (3 runs):
#   run1: "R + 3D, con una Superficie de "
#   run2: "57,40"
#   run3: " m2., en el inmueble individualizado como: Manzana"
#
# Editing run1/run4 in place normally causes the engine to coalesce them
# with their (identically formatted) neighbouring runs -- "...Tipología E"
# before, and the trailing " " run after the next "Manzana" before "${" --
# which would corrupt unrelated text. Toggling Bold on/off around each
# edit is enough to stop that coalescing from reaching outside the
# fragment we actually want to touch, while still ending up with plain
# (non-bold) runs once we clear the flag again.

$d = $word.ActiveDocument

# --- Step 1: trim the trailing "60" off the first run -----------------
$r1 = $d.Content
$r1.Find.Execute("R + 3D, con una Superficie de 60", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r1.Bold = 1
$r1.Text = "R + 3D, con una Superficie de "
$r1Start = $r1.Start
$r1End = $r1.End

# --- Step 2: turn the middle ",3" into the new "57,40" value ----------
$r2 = $d.Content
$r2.Find.Execute(",3", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Bold = 1
$r2.Text = "57,40"
$r2Start = $r2.Start
$r2End = $r2.End

# --- Step 3: trim the leading "0" off the last run, keep the space ----
$r4 = $d.Content
$r4.Find.Execute("0 m2., en el inmueble individualizado como: Manzana", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r4.Bold = 1
$r4.Text = " m2., en el inmueble individualizado como: Manzana"
$r4Start = $r4.Start
$r4End = $r4.End

# --- Clear the temporary Bold shield on each piece individually -------
# (done as three separate Range operations so the now-identical
# formatting doesn't get merged back into one run across segments)
$u1 = $d.Range($r1Start, $r1End)
$u1.Bold = 0
$u2 = $d.Range($r2Start, $r2End)
$u2.Bold = 0
$u4 = $d.Range($r4Start, $r4End)
$u4.Bold = 0
